$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last two data rows (old rows 8 and 9, Resolving-Mac target pairs)
$ws.Rows("8:9").Delete()

# Update remaining data rows (2-7) with the recalculated TPM-based values

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Artn"
$ws.Range("C2").Value = "Ret"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.306751666666667
$ws.Range("H2").Value = 6.920255
$ws.Range("I2").Value = 0.8222690533928816
$ws.Range("J2").Value = 0.8222690533928814
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.132875333333334
$ws.Range("N2").Value = 12.398626
$ws.Range("O2").Value = 0.391502049440379
$ws.Range("P2").Value = 0.3915020494403789
$ws.Range("Q2").Value = 9.533517063292223
$ws.Range("R2").Value = 85.80165356963
$ws.Range("S2").Value = 0.3219200195947136
$ws.Range("T2").Value = 0.3219200195947134

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Artn"
$ws.Range("C3").Value = "Ret"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.306751666666667
$ws.Range("H3").Value = 6.920255
$ws.Range("I3").Value = 0.8222690533928816
$ws.Range("J3").Value = 0.8222690533928814
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.094146333333334
$ws.Range("N3").Value = 15.282439
$ws.Range("O3").Value = 0.4825620346115429
$ws.Range("P3").Value = 0.4825620346115428
$ws.Range("Q3").Value = 11.75093054466056
$ws.Range("R3").Value = 105.758374901945
$ws.Range("S3").Value = 0.3967958274033763
$ws.Range("T3").Value = 0.3967958274033762

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Artn"
$ws.Range("C4").Value = "Ret"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.306751666666667
$ws.Range("H4").Value = 6.920255
$ws.Range("I4").Value = 0.8222690533928816
$ws.Range("J4").Value = 0.8222690533928814
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.329437333333333
$ws.Range("N4").Value = 3.988312
$ws.Range("O4").Value = 0.1259359159480782
$ws.Range("P4").Value = 0.1259359159480782
$ws.Range("Q4").Value = 3.066681784395556
$ws.Range("R4").Value = 27.60013605956
$ws.Range("S4").Value = 0.1035532063947917
$ws.Range("T4").Value = 0.1035532063947917

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Artn"
$ws.Range("C5").Value = "Ret"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.4985973333333333
$ws.Range("H5").Value = 1.495792
$ws.Range("I5").Value = 0.1777309466071186
$ws.Range("J5").Value = 0.1777309466071185
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.132875333333334
$ws.Range("N5").Value = 12.398626
$ws.Range("O5").Value = 0.391502049440379
$ws.Range("P5").Value = 0.3915020494403789
$ws.Range("Q5").Value = 2.060640620199111
$ws.Range("R5").Value = 18.545765581792
$ws.Range("S5").Value = 0.06958202984566549
$ws.Range("T5").Value = 0.06958202984566544

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Artn"
$ws.Range("C6").Value = "Ret"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.4985973333333333
$ws.Range("H6").Value = 1.495792
$ws.Range("I6").Value = 0.1777309466071186
$ws.Range("J6").Value = 0.1777309466071185
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 5.094146333333334
$ws.Range("N6").Value = 15.282439
$ws.Range("O6").Value = 0.4825620346115429
$ws.Range("P6").Value = 0.4825620346115428
$ws.Range("Q6").Value = 2.539927777409778
$ws.Range("R6").Value = 22.859349996688
$ws.Range("S6").Value = 0.08576620720816662
$ws.Range("T6").Value = 0.08576620720816659

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Artn"
$ws.Range("C7").Value = "Ret"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.4985973333333333
$ws.Range("H7").Value = 1.495792
$ws.Range("I7").Value = 0.1777309466071186
$ws.Range("J7").Value = 0.1777309466071185
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.329437333333333
$ws.Range("N7").Value = 3.988312
$ws.Range("O7").Value = 0.1259359159480782
$ws.Range("P7").Value = 0.1259359159480782
$ws.Range("Q7").Value = 0.6628539092337778
$ws.Range("R7").Value = 5.965685183104
$ws.Range("S7").Value = 0.02238270955328645
$ws.Range("T7").Value = 0.02238270955328645
